$wb = $excel.ActiveWorkbook

$wsStats   = $wb.Worksheets.Item("stats")
$wsRace    = $wb.Worksheets.Item("race")
$wsCombo64 = $wb.Worksheets.Item("combo_64")

# ---------------------------------------------------------------------------
# stats sheet: two new rows of notes at the bottom of the table
# ---------------------------------------------------------------------------
$wsStats.Range("B18").Copy() | Out-Null
$wsStats.Range("B20").PasteSpecial(-4122) | Out-Null
# Type B21 first so it claims the shared-string slot before B20's text does.
$wsStats.Range("B21").Value = "формула для пробивания брони простая = def - def(из комбо) - атака"
$wsStats.Range("B20").Value = "формулы для урона = def - атака"

# ---------------------------------------------------------------------------
# race sheet: new ATK column (S) + combo-formula bonus folded into the J col
# ---------------------------------------------------------------------------
$wsRace.Range("J14").Copy() | Out-Null
$wsRace.Range("S2:S11").PasteSpecial(-4122) | Out-Null

$wsRace.Range("S2").Value = "ATK"
$wsRace.Range("S3").Value = 5
$wsRace.Range("S4").Value = 4
$wsRace.Range("S5").Value = 4
$wsRace.Range("S6").Value = 6
$wsRace.Range("S7").Value = 6
$wsRace.Range("S8").Value = 6
$wsRace.Range("S9").Value = 7
$wsRace.Range("S10").Value = 4
$wsRace.Range("S11").Value = 5

$wsRace.Columns.Item(19).ColumnWidth = 3.5

$wsRace.Range("J15").Formula = "=S3+B15/2"
$wsRace.Range("J16:J23").Formula = "=S4+B16/2"

# ---------------------------------------------------------------------------
# view state: selections + which tab/sheet is active
# ---------------------------------------------------------------------------
$wsStats.Range("D13").Select() | Out-Null
$wsCombo64.Range("B21").Select() | Out-Null
$wsRace.Range("U17").Select() | Out-Null
$wsRace.Activate() | Out-Null
